$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.350.57'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').Value = '1.671.08'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('D5').Value = '220.65'
$ws.Range('E5').Value = '  +0.88%  '
$ws.Range('D6').Value = '0.5314'
$ws.Range('E6').Value = '  -0.36%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '0.2656'
$ws.Range('E8').Value = '  +0.56%  '
$ws.Range('D9').Value = '0.06374'
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('D10').Value = '21.02'
$ws.Range('E10').Value = '  +2.29%  '
$ws.Range('D11').Value = '0.07851'
$ws.Range('E11').Value = '  +0.06%  '
$ws.Range('D12').Value = '4.543'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').Value = '1.670.31'
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('D14').Value = '1.898.85'
$ws.Range('E14').Value = '  +0.38%  '
$ws.Range('D15').Value = '0.5634'
$ws.Range('E15').Value = '  +1.75%  '
$ws.Range('D16').Value = '0.0₅8135'
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('D17').Value = '66.01'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('D18').Value = '26.350.56'
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.730'
$ws.Range('E20').Value = '  +1.33%  '
$ws.Range('D21').Value = '201.85'
$ws.Range('E21').Value = '  +4.92%  '
$ws.Range('E22').Value = '  +1.58%  '
$ws.Range('D23').Value = '6.068'
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').Value = '146.66'
$ws.Range('E25').Value = '  +1.06%  '
$ws.Range('D26').Value = '0.1216'
$ws.Range('E26').Value = '  -0.85%  '
$ws.Range('D27').Value = '7.261'
$ws.Range('E27').Value = '  +0.27%  '
$ws.Range('D28').Value = '16.25'
$ws.Range('E28').Value = '  +0.92%  '
$ws.Range('D29').Value = '1.517'
$ws.Range('E29').Value = '  +2.89%  '
$ws.Range('D30').Value = '0.05899'
$ws.Range('E30').Value = '  +1.09%  '
$ws.Range('D31').Value = '1.288'
$ws.Range('E31').Value = '  +0.68%  '
$ws.Range('D32').Value = '3.538'
$ws.Range('E32').Value = '  -1.19%  '
$ws.Range('D33').Value = '3.332'
$ws.Range('E33').Value = '  +0.96%  '
$ws.Range('D34').Value = '1.607'
$ws.Range('D35').Value = '0.9684'
$ws.Range('E35').Value = '  +0.92%  '
$ws.Range('E36').Value = '  +0.31%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').Value = '0.5809'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01620'
$ws.Range('E39').Value = '  +0.76%  '
$ws.Range('D40').Value = '5.979'
$ws.Range('E40').Value = '  +1.61%  '
$ws.Range('D41').Value = '1.079.19'
$ws.Range('E41').Value = '  +3.07%  '
$ws.Range('D42').Value = '0.8602'
$ws.Range('E42').Value = '  +0.66%  '
$ws.Range('D44').Value = '103.18'
$ws.Range('E44').Value = '  -1.66%  '
$ws.Range('D45').Value = '1.808.10'
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('D46').Value = '58.59'
$ws.Range('E46').Value = '  +2.21%  '
$ws.Range('E47').Value = '  +1.14%  '
$ws.Range('D48').Value = '1.017'
$ws.Range('E48').Value = '  +0.83%  '
$ws.Range('E49').Value = '  +1.00%  '
$ws.Range('D50').Value = '8.121'
$ws.Range('E50').Value = '  +2.25%  '
$ws.Range('D51').Value = '0.05149'
$ws.Range('E51').Value = '  -0.33%  '

Write-Host "Updated cryptos list"
